$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.916.17'
$ws.Range('E2').Value = '  +3.55%  '
$ws.Range('D3').Value = '2.700.18'
$ws.Range('E3').Value = '  +1.84%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.45'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +2.21%  '
$ws.Range('D9').Value = '2.722.04'
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('E10').Value = '  +3.65%  '
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('D14').Value = '3.196.28'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').Value = '60.717.51'
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.44'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000139'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.716.98'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '353.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.16%  '
$ws.Range('E22').Value = '  +4.00%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.47%  '
$ws.Range('E25').Value = '  +1.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.169'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').Value = '0.0₃0824'
$ws.Range('E28').Value = '  +2.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.89'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  +1.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '147.74'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.33'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.00%  '
$ws.Range('E36').Value = '  +9.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.958'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.73%  '
$ws.Range('E38').Value = '  +10.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '286.14'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.79%  '
$ws.Range('E44').Value = '  -0.70%  '
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.149.05'
$ws.Range('E46').Value = '  +7.59%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.997'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.93'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.63%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0541'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.50'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.22%  '
